$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 150.26666
$ws.Range("I11").Value = 150.26666
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 150.26666
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -10.26666

$ws.Range("H76").Value = 9066.5
$ws.Range("I76").Value = 9879.799999999999
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 9879.799999999999
$ws.Range("L76").Value = 5000
$ws.Range("M76").Value = -9564.799999999999

$ws.Range("H79").Value = 9066.5
$ws.Range("I79").Value = 9879.799999999999
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 9879.799999999999
$ws.Range("L79").Value = 5000
$ws.Range("M79").Value = -8787.799999999999

$ws.Range("H112").Value = 5910.1567
$ws.Range("I112").Value = 1500
$ws.Range("J112").Value = 5998.36
$ws.Range("K112").Value = 4500
$ws.Range("L112").Value = 17995.08
$ws.Range("M112").Value = -3392
$ws.Range("N112").Value = -20211.08

$ws.Range("H113").Value = 18713.438
$ws.Range("I113").Value = 15520.091
$ws.Range("J113").Value = 25738.8
$ws.Range("K113").Value = 15520.091
$ws.Range("L113").Value = 25738.8
$ws.Range("M113").Value = -12266.091
$ws.Range("N113").Value = -32246.8

$ws.Range("H137").Value = 2418.2456
$ws.Range("I137").Value = 2384.6553
$ws.Range("J137").Value = 2453.0356
$ws.Range("K137").Value = 7153.965899999999
$ws.Range("L137").Value = 7359.1068
$ws.Range("M137").Value = -4603.965899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2185
$ws.Range("I2").Value = 2191.1428
$ws.Range("J2").Value = 2159.2
$ws.Range("K2").Value = 2191.1428
$ws.Range("L2").Value = 2159.2
$ws.Range("M2").Value = -2078.1428

$ws.Range("H45").Value = 7824.75
$ws.Range("I45").Value = 26823.625
$ws.Range("J45").Value = 4024.975
$ws.Range("K45").Value = 26823.625
$ws.Range("L45").Value = 4024.975
$ws.Range("M45").Value = -26446.625

$ws.Range("H62").Value = 59000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 59000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 59000
$ws.Range("N62").Value = -60248

$ws.Range("H65").Value = 59000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 59000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 177000
$ws.Range("N65").Value = -183240

$ws.Range("H110").Value = 3209.348
$ws.Range("I110").Value = 3242.2354
$ws.Range("J110").Value = 3116.1667
$ws.Range("K110").Value = 3242.2354
$ws.Range("L110").Value = 3116.1667
$ws.Range("M110").Value = -1197.2354

$ws.Range("H116").Value = 2185
$ws.Range("I116").Value = 2191.1428
$ws.Range("J116").Value = 2159.2
$ws.Range("K116").Value = 2191.1428
$ws.Range("L116").Value = 2159.2
$ws.Range("M116").Value = 102.8571999999999

$ws.Range("H132").Value = 2218.4
$ws.Range("I132").Value = 1989.5294
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 5968.5882
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -3438.5882

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2185
$ws.Range("I3").Value = 2191.1428
$ws.Range("J3").Value = 2159.2
$ws.Range("K3").Value = 2191.1428
$ws.Range("L3").Value = 2159.2
$ws.Range("M3").Value = -2077.1428

$ws.Range("H86").Value = 15743326
$ws.Range("I86").Value = 23613536
$ws.Range("J86").Value = 2907.3333
$ws.Range("K86").Value = 23613536
$ws.Range("L86").Value = 2907.3333
$ws.Range("M86").Value = -23612413

$ws.Range("H88").Value = 20998.334

$ws.Range("H89").Value = 15743326
$ws.Range("I89").Value = 23613536
$ws.Range("J89").Value = 2907.3333
$ws.Range("K89").Value = 118067680
$ws.Range("L89").Value = 14536.6665
$ws.Range("M89").Value = -118062064

$ws.Range("H91").Value = 20998.334

$ws.Range("H105").Value = 2979.2
$ws.Range("I105").Value = 2865.7778
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 2865.7778
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -1118.7778

$ws.Range("H107").Value = 3817.9614
$ws.Range("I107").Value = 3865.375
$ws.Range("J107").Value = 3249
$ws.Range("K107").Value = 3865.375
$ws.Range("L107").Value = 3249
$ws.Range("M107").Value = -1945.375

$ws.Range("H134").Value = 2722.9688
$ws.Range("I134").Value = 2714.0645
$ws.Range("J134").Value = 2999
$ws.Range("K134").Value = 8142.193499999999
$ws.Range("L134").Value = 8997
$ws.Range("M134").Value = -5607.193499999999
$ws.Range("N134").Value = -14067

$ws.Range("H138").Value = 88774.5
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 88774.5
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 88774.5
$ws.Range("N138").Value = -99054.5

$ws.Range("H140").Value = 88000
$ws.Range("I140").Value = 88000
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 88000
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -82820
$ws.Range("N140").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3792.0435
$ws.Range("I58").Value = 4097.8613
$ws.Range("J58").Value = 2691.1
$ws.Range("K58").Value = 4097.8613
$ws.Range("L58").Value = 2691.1
$ws.Range("M58").Value = -3894.8613

$ws.Range("H62").Value = 2734.2222
$ws.Range("I62").Value = 2734.2222
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2734.2222
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2110.2222

$ws.Range("H65").Value = 2734.2222
$ws.Range("I65").Value = 2734.2222
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 13671.111
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -10551.111

$ws.Range("H136").Value = 3792.0435
$ws.Range("I136").Value = 4097.8613
$ws.Range("J136").Value = 2691.1
$ws.Range("K136").Value = 12293.5839
$ws.Range("L136").Value = 8073.299999999999
$ws.Range("M136").Value = -9743.583899999998

$ws.Range("H138").Value = 92656.73
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 92656.73
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 92656.73
$ws.Range("N138").Value = -102936.73

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = $null
$ws.Range("N81").Value = $null

$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = $null
$ws.Range("N84").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 16259
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 16259
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 16259
$ws.Range("N46").Value = -16571

$ws.Range("H58").Value = 21823.8
$ws.Range("I58").Value = 21279.75
$ws.Range("J58").Value = 24000
$ws.Range("K58").Value = 21279.75
$ws.Range("L58").Value = 24000
$ws.Range("M58").Value = -21002.75
$ws.Range("N58").Value = -24554

$ws.Range("H80").Value = 8441
$ws.Range("I80").Value = 5999.5
$ws.Range("J80").Value = 9417.6
$ws.Range("K80").Value = 5999.5
$ws.Range("L80").Value = 9417.6
$ws.Range("M80").Value = -5001.5
$ws.Range("N80").Value = -11413.6

$ws.Range("H83").Value = 8441
$ws.Range("I83").Value = 5999.5
$ws.Range("J83").Value = 9417.6
$ws.Range("K83").Value = 29997.5
$ws.Range("L83").Value = 47088
$ws.Range("M83").Value = -25005.5
$ws.Range("N83").Value = -57072

$ws.Range("H132").Value = 4476.5317
$ws.Range("I132").Value = 4582.943
$ws.Range("J132").Value = 4166.1665
$ws.Range("K132").Value = 13748.829
$ws.Range("L132").Value = 12498.4995
$ws.Range("M132").Value = -11218.829
$ws.Range("N132").Value = -17558.4995

$ws.Range("H135").Value = 84886.75
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 84886.75
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 84886.75
$ws.Range("N135").Value = -95026.75

$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H57").Value = 22500
$ws.Range("I57").Value = 20000
$ws.Range("J57").Value = 25000
$ws.Range("K57").Value = 20000
$ws.Range("L57").Value = 25000
$ws.Range("M57").Value = -19434
$ws.Range("N57").Value = -26132

$ws.Range("H62").Value = 54997
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 54997
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 54997
$ws.Range("N62").Value = -56245

$ws.Range("H65").Value = 54997
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 54997
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 164991
$ws.Range("N65").Value = -171231

$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = $null

$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = $null

$ws.Range("H137").Value = 89987
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 89987
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 89987
$ws.Range("N137").Value = -100187

$ws.Range("H139").Value = 49999
$ws.Range("I139").Value = 49999
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 49999
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -44859
$ws.Range("N139").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3618.9211
$ws.Range("I122").Value = 3302.3
$ws.Range("J122").Value = 4806.25
$ws.Range("K122").Value = 9906.900000000001
$ws.Range("L122").Value = 14418.75
$ws.Range("M122").Value = -7456.900000000001
$ws.Range("N122").Value = -19318.75

$ws.Range("H132").Value = 4151.851
$ws.Range("I132").Value = 3856.25
$ws.Range("J132").Value = 5841
$ws.Range("K132").Value = 11568.75
$ws.Range("L132").Value = 17523
$ws.Range("M132").Value = -9038.75
$ws.Range("N132").Value = -22583
